$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename headers
$ws.Range("A1").Value = "Sample Years"
$ws.Range("E1").Value = "Output Years"

# Update input values (fewer optimization runs, starting at a no-outage point)
$ws.Range("F2").Value = 9
$ws.Range("G2").Value = 700

# Adjust column A width to fit the new, longer header text
$ws.Columns.Item(1).ColumnWidth = 10.510416666666666
